$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only the date changes
$ws.Range("D2").Value = 44284

# Row 4: takes on the values that used to be in row 5
$ws.Range("D4").Value = 44291
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 550

# Row 5: takes on the values that used to be in row 2
$ws.Range("D5").Value = 44280
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 500
